# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# ----- Rushing sheet -----
$rushing = $wb.Worksheets.Item("Rushing")

# S.Ehlinger (row 2)
$rushing.Range("D2").Value = 9
$rushing.Range("E2").Value = 18

# N.Hines (row 4)
$rushing.Range("C4").Value = 147
$rushing.Range("D4").Value = 95
$rushing.Range("E4").Value = 28
$rushing.Range("F4").Value = 75

# A.Dulin (row 10)
$rushing.Range("C10").Value = 3

# ----- Receiving sheet -----
$receiving = $wb.Worksheets.Item("Receiving")

# J.Taylor (row 2)
$receiving.Range("C2").Value = 41
$receiving.Range("D2").Value = 32
$receiving.Range("G2").Value = 5
$receiving.Range("H2").Value = 4

# N.Hines (row 3)
$receiving.Range("C3").Value = 43

# Z.Pascal (row 5)
$receiving.Range("C5").Value = 81
$receiving.Range("D5").Value = 58
$receiving.Range("E5").Value = 22
$receiving.Range("G5").Value = 14

# A.Dulin (row 6)
$receiving.Range("C6").Value = 52
$receiving.Range("E6").Value = 11
$receiving.Range("F6").Value = 5

# D.Patmon (row 10)
$receiving.Range("C10").Value = 33
$receiving.Range("D10").Value = 27

# J.Doyle (row 13)
$receiving.Range("C13").Value = 37
$receiving.Range("D13").Value = 26
$receiving.Range("G13").Value = 8
$receiving.Range("H13").Value = 5

# M.Alie-Cox (row 14)
$receiving.Range("C14").Value = 23
$receiving.Range("G14").Value = 8
